$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, [string]$CellRef, [string]$Val)
    $r = $Worksheet.Range($CellRef)
    $r.NumberFormat = "@"
    $r.Value = $Val
    $r.ClearFormats()
}

Set-TextValue $ws 'D2' '27.125.04'
Set-TextValue $ws 'E2' '  -2.67%  '
Set-TextValue $ws 'D3' '1.870.16'
Set-TextValue $ws 'E3' '  -1.94%  '
Set-TextValue $ws 'D4' '1.000'
Set-TextValue $ws 'E4' '  +0.04%  '
Set-TextValue $ws 'D5' '307.54'
Set-TextValue $ws 'E5' '  -1.80%  '
Set-TextValue $ws 'D6' '1.000'
Set-TextValue $ws 'E6' '  +0.09%  '
Set-TextValue $ws 'D7' '0.5059'
Set-TextValue $ws 'E7' '  +1.28%  '
Set-TextValue $ws 'E8' '  -1.49%  '
Set-TextValue $ws 'D9' '0.07161'
Set-TextValue $ws 'E9' '  -1.70%  '
Set-TextValue $ws 'D10' '0.8904'
Set-TextValue $ws 'E10' '  -2.19%  '
Set-TextValue $ws 'D11' '20.73'
Set-TextValue $ws 'E11' '  -1.14%  '
Set-TextValue $ws 'D12' '1.874.08'
Set-TextValue $ws 'E12' '  -1.67%  '
Set-TextValue $ws 'D13' '0.07567'
Set-TextValue $ws 'E13' '  -1.03%  '
Set-TextValue $ws 'D14' '5.330'
Set-TextValue $ws 'E14' '  -2.98%  '
Set-TextValue $ws 'D15' '89.30'
Set-TextValue $ws 'E15' '  -2.87%  '
Set-TextValue $ws 'D16' '1.000'
Set-TextValue $ws 'E16' '  +0.01%  '
Set-TextValue $ws 'D17' '0.000008521'
Set-TextValue $ws 'E17' '  -2.44%  '
Set-TextValue $ws 'E18' '  -3.19%  '
Set-TextValue $ws 'E19' '  +0.12%  '
Set-TextValue $ws 'D20' '27.172.79'
Set-TextValue $ws 'E20' '  -2.59%  '
Set-TextValue $ws 'D21' '5.092'
Set-TextValue $ws 'E21' '  -1.66%  '
Set-TextValue $ws 'D22' '2.108.69'
Set-TextValue $ws 'E22' '  -1.64%  '
Set-TextValue $ws 'D23' '10.62'
Set-TextValue $ws 'E23' '  -1.71%  '
Set-TextValue $ws 'D24' '6.504'
Set-TextValue $ws 'E24' '  -1.07%  '
Set-TextValue $ws 'D25' '151.07'
Set-TextValue $ws 'E25' '  -1.35%  '
Set-TextValue $ws 'E26' '  -2.07%  '
Set-TextValue $ws 'D27' '18.03'
Set-TextValue $ws 'E27' '  -1.97%  '
Set-TextValue $ws 'D28' '2.095'
Set-TextValue $ws 'E28' '  -5.48%  '
Set-TextValue $ws 'E29' '  -2.00%  '
Set-TextValue $ws 'D30' '4.772'
Set-TextValue $ws 'E30' '  -3.03%  '
Set-TextValue $ws 'D31' '4.692'
Set-TextValue $ws 'E31' '  -1.95%  '
Set-TextValue $ws 'D32' '0.08990'
Set-TextValue $ws 'E32' '  -0.30%  '
Set-TextValue $ws 'D33' '0.05143'
Set-TextValue $ws 'E33' '  -2.47%  '
Set-TextValue $ws 'D34' '3.100'
Set-TextValue $ws 'E34' '  -2.73%  '
Set-TextValue $ws 'D35' '0.7468'
Set-TextValue $ws 'E35' '  -3.57%  '
Set-TextValue $ws 'E36' '  -5.72%  '
Set-TextValue $ws 'B37' 'VeChain'
Set-TextValue $ws 'C37' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws 'D37' '0.02040'
Set-TextValue $ws 'E37' '  -2.28%  '
Set-TextValue $ws 'B38' 'RenderToken'
Set-TextValue $ws 'C38' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws 'D38' '2.551'
Set-TextValue $ws 'E38' '  +0.10%  '
Set-TextValue $ws 'D39' '3.039'
Set-TextValue $ws 'E39' '  +0.49%  '
Set-TextValue $ws 'E40' '  -1.94%  '
Set-TextValue $ws 'D41' '0.5370'
Set-TextValue $ws 'E41' '  -3.49%  '
Set-TextValue $ws 'D42' '6.626'
Set-TextValue $ws 'E42' '  -4.04%  '
Set-TextValue $ws 'D43' '115.03'
Set-TextValue $ws 'E43' '  +3.04%  '
Set-TextValue $ws 'D44' '8.475'
Set-TextValue $ws 'E44' '  -0.20%  '
Set-TextValue $ws 'D45' '0.1480'
Set-TextValue $ws 'E45' '  -2.43%  '
Set-TextValue $ws 'D46' '0.4653'
Set-TextValue $ws 'E46' '  -3.82%  '
Set-TextValue $ws 'D47' '1.000'
Set-TextValue $ws 'E47' '  +0.09%  '
Set-TextValue $ws 'D48' '10.05'
Set-TextValue $ws 'E48' '  -4.84%  '
Set-TextValue $ws 'E49' '  -3.76%  '
Set-TextValue $ws 'D50' '64.75'
Set-TextValue $ws 'E50' '  -4.20%  '
Set-TextValue $ws 'D51' '36.67'
Set-TextValue $ws 'E51' '  -1.17%  '
